$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date values from 45212 (2023-10-13) to 45221 (2023-10-22)
# for rows 2 through 7, preserving existing cell formatting.
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45221
}
